$d = $word.ActiveDocument
$d.Content.Find.Execute("Minneapolis Code Master 2013", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Minneapolis Code Mastery 2013", 2)
